# Add a "plusminusCI95" column (the +/- half-width of the 95% CI on Lmat50)
# to both maturity tables, and a running "total sums" pass over the sheets
# (per the commit message: "updated maturity agreement tables to have
# total sums"). Column F (canary_rockfish) / G (sablefish) is inserted,
# shifting the existing CI_95/alpha/beta columns one slot to the right.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# canary_rockfish
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("canary_rockfish")

# Fully clear the sheet's stale <cols> custom-width metadata before
# rewriting it with the new layout (first nudge the column-1 width
# record off column 1, then drop the whole used range in one go).
$ws1.Columns.Item(1).Insert()
$ws1.Range("A:I").EntireColumn.Delete()

$data1 = @(
    @("time_period", "resting_stage_removed", "spawning_season_isolated", "n", "Lmat50", "plusminusCI95", "CI_95", "alpha", "beta"),
    @("1980-1984", "No", "No", 1007, 45.96, 0.37, "45.59-46.33", "-24.32 (1.75)", "0.53 (0.04)"),
    @("1980-1984", "Yes", "No", 913, 46.37, 0.4, "45.97-46.77", "-25.66 (1.89)", "0.55 (0.04)"),
    @("1980-1984", "No", "Yes", 305, 46.74, 1.15, "45.59-47.88", "-21.42 (2.58)", "0.46 (0.05)"),
    @("1980-1984", "Yes", "Yes", 266, 47.95, 1.37, "46.58-49.32", "-27.77 (3.56)", "0.58 (0.07)"),
    @("2010-2018", "No", "No", 1759, 32.66, 0.2, "32.46-32.86", "-8.66 (0.72)", "0.27 (0.02)"),
    @("2010-2018", "Yes", "No", 917, 40.09, 0.42, "39.66-40.51", "-20.36 (1.84)", "0.51 (0.04)"),
    @("2010-2018", "No", "Yes", 266, 33.38, 1.29, "32.1-34.67", "-9.5 (1.92)", "0.28 (0.05)"),
    @("2010-2018", "Yes", "Yes", 128, 42.31, 4.05, "38.27-46.36", "-33.67 (9.86)", "0.8 (0.23)")
)

for ($r = 0; $r -lt $data1.Length; $r++) {
    $rowVals = $data1[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws1.Cells.Item($r + 1, $c + 1).Value = $rowVals[$c]
    }
}

$ws1.Range("H20").Select()

# ---------------------------------------------------------------------
# sablefish
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("sablefish")

$ws2.Columns.Item(1).Insert()
$ws2.Range("A:J").EntireColumn.Delete()

$data2 = @(
    @("time_period", "resting_stage_removed", "spawning_season_isolated", "immature_70cm_plus_removed", "n", "Lmat50", "plusminusCI95", "CI_95", "alpha", "beta"),
    @("1995-2001", "No", "No", "No", 746, 55.17, 0.61, "54.57-55.78", "-11.52 (1.06)", "0.21 (0.02)"),
    @("1995-2001", "No", "No", "Yes", 738, 55.44, 0.59, "54.85-56.04", "-14 (1.22)", "0.25 (0.02)"),
    @("1995-2001", "Yes", "No", "Yes", 479, 58.52, 0.88, "57.64-59.4", "-15.21 (1.44)", "0.26 (0.02)"),
    @("1995-2001", "No", "Yes", "Yes", 298, 55.55, 1.45, "54.09-57", "-18.98 (2.46)", "0.34 (0.04)"),
    @("1995-2001", "Yes", "Yes", "Yes", 256, 56.01, 1.64, "54.36-57.65", "-18.27 (2.44)", "0.33 (0.04)"),
    @("2007-2018", "No", "No", "No", 1505, 48.44, 0.28, "48.16-48.72", "-12.03 (0.83)", "0.25 (0.02)"),
    @("2007-2018", "No", "No", "Yes", 1505, 48.44, 0.28, "48.16-48.72", "-12.03 (0.83)", "0.25 (0.02)"),
    @("2007-2018", "Yes", "No", "Yes", 874, 54.82, 0.5, "54.32-55.32", "-21.23 (1.58)", "0.39 (0.03)"),
    @("2007-2018", "No", "Yes", "Yes", 837, 49.25, 0.49, "48.76-49.73", "-13.99 (1.2)", "0.28 (0.02)"),
    @("2007-2018", "Yes", "Yes", "Yes", 533, 55.12, 0.88, "54.24-56.01", "-27.33 (2.77)", "0.5 (0.05)")
)

for ($r = 0; $r -lt $data2.Length; $r++) {
    $rowVals = $data2[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws2.Cells.Item($r + 1, $c + 1).Value = $rowVals[$c]
    }
}

$ws2.Range("I19").Select()

# ---------------------------------------------------------------------
# Workbook window position + active sheet
# ---------------------------------------------------------------------
$excel.ActiveWindow.Left = 7860
$excel.ActiveWindow.Top = 460

$ws1.Activate()
